# Apply updated dSF (column F) values per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -3
$ws.Range("F4").Value = -8
$ws.Range("F5").Value = -7
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = 10
$ws.Range("F9").Value = -10
$ws.Range("F12").Value = -1
$ws.Range("F13").Value = -1
$ws.Range("F14").Value = 9
$ws.Range("F17").Value = 4
$ws.Range("F18").Value = -1
